$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (ASSISTS) on rows 2-41 is currently stored as text ("0", "1", ...).
# Re-assign as real numbers so the cells become numeric (t="n") instead of inlineStr.
$ws.Range("F2:F17").Value = 0
$ws.Range("F18:F36").Value = 1
$ws.Range("F37:F38").Value = 2
$ws.Range("F39:F41").Value = 3

# Column H (CHAMPION) fixes: several rows were mislabeled and should read "Olaf".
$ws.Range("H5").Value = "Olaf"
$ws.Range("H11").Value = "Olaf"
$ws.Range("H15").Value = "Olaf"
$ws.Range("H17").Value = "Olaf"
$ws.Range("H21").Value = "Olaf"
$ws.Range("H23").Value = "Olaf"
$ws.Range("H29").Value = "Olaf"
$ws.Range("H35").Value = "Olaf"
$ws.Range("H41").Value = "Olaf"
